$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 20): two parallel mini-tables (A:F and G:L) ---
$ws.Range("A20").Value = "x"
$ws.Range("A20").Style = "Input"
$ws.Range("B20").Value = "y distribution type"
$ws.Range("B20").Style = "Neutral"
$ws.Range("C20").Value = "DistVal1"
$ws.Range("D20").Value = "DistVal2"
$ws.Range("E20").Value = "DistVal3"
$ws.Range("F20").Value = "DistVal4"

$ws.Range("G20").Value = "x"
$ws.Range("G20").Style = "Input"
$ws.Range("H20").Value = "y distribution type"
$ws.Range("H20").Style = "Neutral"
$ws.Range("I20").Value = "DistVal1"
$ws.Range("J20").Value = "DistVal2"
$ws.Range("K20").Value = "DistVal3"
$ws.Range("L20").Value = "DistVal4"

# --- Row 21: first data row with sample distribution values ---
$ws.Range("A21").Value = 1
$ws.Range("B21").Value = "Normal"
$ws.Range("C21").Value = 100
$ws.Range("D21").Value = 1

$ws.Range("G21").Value = 1
$ws.Range("H21").Value = "Triangular"
$ws.Range("I21").Value = 100
$ws.Range("J21").Value = 150
$ws.Range("K21").Value = 200

# --- Rows 22-25: remaining index values for both tables ---
$ws.Range("A22").Value = 2
$ws.Range("G22").Value = 2

$ws.Range("A23").Value = 3
$ws.Range("G23").Value = 3

$ws.Range("A24").Value = 4
$ws.Range("G24").Value = 4

$ws.Range("A25").Value = 5
$ws.Range("G25").Value = 5

# --- Column H width ---
$ws.Columns.Item(8).ColumnWidth = 10.736979166666666

# --- Selection matching the recorded cursor position ---
$ws.Range("B20:F22").Select()
